# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Update the time_taken (F column) timestamps on the "data" sheet ---
$timestamps = @(
    "2021-10-05 14:21:08.609839",
    "2021-10-05 14:21:08.609846",
    "2021-10-05 14:21:08.609850",
    "2021-10-05 14:21:08.609852",
    "2021-10-05 14:21:08.609855",
    "2021-10-05 14:21:08.609858",
    "2021-10-05 14:21:08.609860",
    "2021-10-05 14:21:08.609863",
    "2021-10-05 14:21:08.609866",
    "2021-10-05 14:21:08.609869",
    "2021-10-05 14:21:08.609871",
    "2021-10-05 14:21:08.609874"
)
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Range("F$row").Value = $timestamps[$i]
}

# --- 2. Add the new "metadata" worksheet, placed after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# --- 3. Header row (bold / centered / bordered, matching the "data" sheet's style) ---
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
$headerCols = @("B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $metaSheet.Range($headerCols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- 4. Data row 2 ---
$a2 = $metaSheet.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$metaSheet.Range("B2").Value = "Inherited pancreatic cancer"
$metaSheet.Range("C2").Value = 524
$metaSheet.Range("D2").Value = "'1.18"
$metaSheet.Range("E2").Value = "2021-04-14T09:20:55.642319Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:21:08.606145"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/524/?format=json"

$dataSheet.Select()
